# Budget sheet update: add a new expense line item for car rentals and
# let the existing SUM formulas (Travel Total / Grand Total) pick it up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New expense row under the "Travel" section (columns A/B).
$ws.Range("A4").Value = "Car Rental (5 Exotic Cars)"
$ws.Range("B4").Value = 59127
$ws.Range("B4").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# Restore the cursor/selection to where it ended up after the edit.
$null = $ws.Range("C8").Select()
